# Helper: build the VBA/COM "Long" RGB value (R + G*256 + B*65536) that the
# PowerPoint object model expects when assigning ColorFormat.RGB.
function RGBVal($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Table on slide 6 switches to a different built-in table style.
# ---------------------------------------------------------------------------
$slide = $p.Slides.Item(6)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{4B3315D6-FE1A-4BD3-AC2F-21FF75E16418}")
    }
}

# ---------------------------------------------------------------------------
# 2) Swap the presentation's theme palette from "Integral" to "Office Theme".
#    (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink - in that order.)
# ---------------------------------------------------------------------------
$theme = $p.SlideMaster.Theme
$theme.Name = "Office Theme"

$colorScheme = $theme.ThemeColorScheme

$colorScheme.Colors(1).RGB  = RGBVal 0x00 0x00 0x00   # dk1
$colorScheme.Colors(2).RGB  = RGBVal 0xFF 0xFF 0xFF   # lt1
$colorScheme.Colors(3).RGB  = RGBVal 0x44 0x54 0x6A   # dk2
$colorScheme.Colors(4).RGB  = RGBVal 0xE7 0xE6 0xE6   # lt2
$colorScheme.Colors(5).RGB  = RGBVal 0x5B 0x9B 0xD5   # accent1
$colorScheme.Colors(6).RGB  = RGBVal 0xED 0x7D 0x31   # accent2
$colorScheme.Colors(7).RGB  = RGBVal 0xA5 0xA5 0xA5   # accent3
$colorScheme.Colors(8).RGB  = RGBVal 0xFF 0xC0 0x00   # accent4
$colorScheme.Colors(9).RGB  = RGBVal 0x44 0x72 0xC4   # accent5
$colorScheme.Colors(10).RGB = RGBVal 0x70 0xAD 0x47   # accent6
$colorScheme.Colors(11).RGB = RGBVal 0x05 0x63 0xC1   # hlink
$colorScheme.Colors(12).RGB = RGBVal 0x95 0x4F 0x72   # folHlink
